$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.733.99"
$ws.Range("E2").Value = "  +8.37%  "

$ws.Range("D3").Value = "3.466.31"
$ws.Range("E3").Value = "  +5.48%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "414.11"
$ws.Range("E5").Value = "  +3.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "124.62"
$ws.Range("E6").Value = "  +13.62%  "

$ws.Range("D7").Value = "3.455.53"
$ws.Range("E7").Value = "  +5.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +1.60%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.683"
$ws.Range("E10").Value = "  +9.34%  "

$ws.Range("E11").Value = "  +32.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.36"
$ws.Range("E12").Value = "  +4.65%  "

$ws.Range("E13").Value = "  +0.46%  "

$ws.Range("D14").Value = "4.017.99"
$ws.Range("E14").Value = "  +5.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.60"
$ws.Range("E15").Value = "  +3.17%  "

$ws.Range("E16").Value = "  +4.62%  "

$ws.Range("D17").Value = "3.461.46"
$ws.Range("E17").Value = "  +5.17%  "

$ws.Range("D18").Value = "62.760.63"
$ws.Range("E18").Value = "  +8.91%  "

$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.81"
$ws.Range("E20").Value = "  -1.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000137"
$ws.Range("E21").Value = "  +27.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.33"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "316.37"
$ws.Range("E23").Value = "  +5.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.11"
$ws.Range("E24").Value = "  +10.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.96"
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.86"
$ws.Range("E27").Value = "  +9.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").Value = "  +5.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("E29").Value = "  -0.83%  "

$ws.Range("E30").Value = "  +3.17%  "

$ws.Range("E31").Value = "  -2.08%  "

$ws.Range("E32").Value = "  +4.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.64"
$ws.Range("E33").Value = "  +23.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.56"
$ws.Range("E34").Value = "  +2.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.11"
$ws.Range("E35").Value = "  +3.28%  "

$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("E37").Value = "  -1.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.26"
$ws.Range("E38").Value = "  +0.94%  "

$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("E41").Value = "  -3.53%  "

$ws.Range("E42").Value = "  +6.11%  "

$ws.Range("E43").Value = "  +2.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.74"
$ws.Range("E44").Value = "  -2.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.285"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.91"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.90"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.25"
$ws.Range("E48").Value = "  +0.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.94"
$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("D50").Value = "2.205.48"
$ws.Range("E50").Value = "  +2.05%  "

$ws.Range("E51").Value = "  +1.15%  "
